$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "MCT-2A-Eletrônica analóg. e de potência"

$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "MCT-2A-Eletrônica analóg. e de potência"
$ws.Range("E4").Value = "-"
$ws.Range("F4").Value = "-"

$ws.Range("F6").Value = "-"

$ws.Range("C8").Value = "-"
$ws.Range("E8").Value = "-"

$ws.Range("B10").Value = "-"
$ws.Range("F10").Value = "[-, 'MEC-2A-Elet. Dig. Bas.', -, -]"

$ws.Range("B11").Value = "[-, 'MEC-1A-Comandos Eletricos', -, -]"
$ws.Range("C11").Value = "-"
$ws.Range("E11").Value = "-"

$ws.Range("B12").Value = "[-, 'MEC-1A-Comandos Eletricos', -, -]"
$ws.Range("C12").Value = "-"
$ws.Range("E12").Value = "-"

$ws.Range("B14").Value = "[-, 'MEC-1A-Comandos Eletricos', -, -]"
$ws.Range("C14").Value = "-"
$ws.Range("E14").Value = "-"
$ws.Range("F14").Value = "[-, -, -, 'MEC-2A-Elet. Dig. Bas.']"

$ws.Range("B15").Value = "[-, 'MEC-1A-Comandos Eletricos', -, -]"
$ws.Range("C15").Value = "-"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "[-, -, -, 'MEC-2A-Elet. Dig. Bas.']"

$ws.Range("B16").Value = "-"
$ws.Range("F16").Value = "[-, -, -, 'MEC-2A-Elet. Dig. Bas.']"

$ws.Range("B18").Value = "-"
$ws.Range("D18").Value = "ELM-2NA-Circuitos Elétricos 2"
$ws.Range("E18").Value = "['MEC-1NB-Elet. Dig. Bas.', -, -, -]"

$ws.Range("B19").Value = "-"
$ws.Range("C19").Value = "-"
$ws.Range("D19").Value = "ELM-2NA-Circuitos Elétricos 2"
$ws.Range("E19").Value = "['MEC-1NB-Elet. Dig. Bas.', -, -, -]"

$ws.Range("B20").Value = "-"
$ws.Range("E20").Value = "['MEC-1NB-Elet. Dig. Bas.', -, -, -]"

$ws.Range("B21").Value = "-"
$ws.Range("E21").Value = "['MEC-1NB-Elet. Dig. Bas.', -, -, -]"
